$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.188.72"
$ws.Range("E2").Value = "  -3.52%  "
$ws.Range("D3").Value = "2.456.85"
$ws.Range("E3").Value = "  -2.78%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "309.58"
$ws.Range("E5").Value = "  +0.19%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "92.76"
$ws.Range("E6").Value = "  -7.83%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.550"
$ws.Range("E7").Value = "  -3.02%  "
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.492"
$ws.Range("E9").Value = "  -5.97%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "33.10"
$ws.Range("E10").Value = "  -7.39%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0773"
$ws.Range("E11").Value = "  -3.90%  "
$ws.Range("E12").Value = "  -1.03%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.93"
$ws.Range("E13").Value = "  -5.34%  "
$ws.Range("D14").Value = "2.838.40"
$ws.Range("E14").Value = "  -2.71%  "
$ws.Range("D15").Value = "2.485.43"
$ws.Range("E15").Value = "  -3.39%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.58"
$ws.Range("E16").Value = "  -5.30%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.777"
$ws.Range("E17").Value = "  -3.93%  "
$ws.Range("D18").Value = "41.192.75"
$ws.Range("E18").Value = "  -3.49%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.27"
$ws.Range("E19").Value = "  -6.72%  "
$ws.Range("D20").Value = "0.0₃0914"
$ws.Range("E20").Value = "  -3.75%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.21"
$ws.Range("E21").Value = "  -8.89%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "67.55"
$ws.Range("E22").Value = "  -2.75%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "234.89"
$ws.Range("E23").Value = "  -3.63%  "
$ws.Range("E24").Value = "  -4.45%  "
$ws.Range("E25").Value = "  +0.19%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.89"
$ws.Range("E26").Value = "  -6.93%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "23.81"
$ws.Range("E27").Value = "  -6.46%  "
$ws.Range("E28").Value = "  -5.63%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.55"
$ws.Range("E29").Value = "  -5.72%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "35.43"
$ws.Range("E30").Value = "  -8.57%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "151.20"
$ws.Range("E31").Value = "  -4.29%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.45"
$ws.Range("E32").Value = "  -4.97%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.65"
$ws.Range("E33").Value = "  -5.73%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.56"
$ws.Range("E34").Value = "  -2.80%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0733"
$ws.Range("E35").Value = "  -6.58%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.96"
$ws.Range("E36").Value = "  -6.15%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.86"
$ws.Range("E37").Value = "  -6.84%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "16.74"
$ws.Range("E38").Value = "  -5.76%  "
$ws.Range("B39").Value = "Stellar"
$ws.Range("C39").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.113"
$ws.Range("E39").Value = "  -3.97%  "
$ws.Range("B40").Value = "Kaspa"
$ws.Range("C40").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.102"
$ws.Range("E40").Value = "  -7.69%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.13"
$ws.Range("E41").Value = "  -1.81%  "
$ws.Range("B42").Value = "FirstDigitalUSD"
$ws.Range("C42").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.00"
$ws.Range("E42").Value = "  +0.17%  "
$ws.Range("B43").Value = "EnergySwap"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "20.13"
$ws.Range("E43").Value = "  -7.97%  "
$ws.Range("D44").Value = "1.980.73"
$ws.Range("E44").Value = "  -1.40%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0282"
$ws.Range("E45").Value = "  -6.38%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.99"
$ws.Range("E46").Value = "  -8.91%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.52"
$ws.Range("E47").Value = "  -4.22%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "70.49"
$ws.Range("E48").Value = "  -2.41%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "95.74"
$ws.Range("E49").Value = "  -5.56%  "
$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.175"
$ws.Range("E50").Value = "  -7.61%  "
$ws.Range("B51").Value = "BitcoinSV"
$ws.Range("C51").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "73.37"
$ws.Range("E51").Value = "  -7.20%  "
